$wb = $excel.ActiveWorkbook

# NOTE: "Vector_bf" and "Vector_BF" are two distinct sheets that differ only
# by case, and Worksheets.Item(name) resolves case-insensitively (both names
# would otherwise hit the same sheet) - so every sheet below is addressed by
# its fixed 1-based position instead of by name.
$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsPunto    = $wb.Worksheets.Item(4)   # Punto_modificado
$wsVecbf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)   # Vector_BF

# Helper: write a numeric-looking value as TEXT (matching the source data,
# which stores generated numbers as plain strings) without leaving a
# lingering custom number-format style on the cell.
function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- Restricciones_del_follower ---------------------------------------
# Row 2 (-y_1)
Set-TextValue $wsFollower "B2" "-4.657691821664619"
Set-TextValue $wsFollower "D2" "0.24011722556595838"
Set-TextValue $wsFollower "E2" "0.24674709760210112"
Set-TextValue $wsFollower "F2" "0.33551389136612164"

# Row 3 (-4 + y_1)
Set-TextValue $wsFollower "B3" "0.6576918216646188"
Set-TextValue $wsFollower "D3" "0.1083236165390392"
Set-TextValue $wsFollower "E3" "0.45750866817633074"
Set-TextValue $wsFollower "F3" "0.88288873043222"

# Row 4 (-16 - 2x + 5y_1)
Set-TextValue $wsFollower "B4" "-4.4622558915346"
Set-TextValue $wsFollower "D4" "0.1102758390135593"
Set-TextValue $wsFollower "E4" "0.9701503029633237"
Set-TextValue $wsFollower "F4" "0.6534924579065517"

# Row 5 (-48 + 8x + y_1)
Set-TextValue $wsFollower "B5" "3.6605518210954013"
Set-TextValue $wsFollower "D5" "0.3168885247170169"
Set-TextValue $wsFollower "E5" "0.5253190516769398"
Set-TextValue $wsFollower "F5" "0.02911925410430937"

# Row 6 (12 - 2x - 2y_1)
Set-TextValue $wsFollower "B6" "-9.066098643186933"
Set-TextValue $wsFollower "D6" "0.4167665579899481"
Set-TextValue $wsFollower "E6" "0.6274738755365885"
Set-TextValue $wsFollower "F6" "0.8208659928897446"

# --- Punto_modificado ----------------------------------------------------
Set-TextValue $wsPunto "A2" "5.875357499928848"
Set-TextValue $wsPunto "B2" "4.657691821664619"
Set-TextValue $wsPunto "C2" "2.112315956957238"

# --- Vector_bf -------------------------------------------------------------
Set-TextValue $wsVecbf "A2" "1.0970590052220022"

# --- Vector_BF -------------------------------------------------------------
Set-TextValue $wsVecBF "A2" "-0.0073040564156943155"
Set-TextValue $wsVecBF "A3" "-1.331884385994611"
